# Update the "Corr/total marks" figures on the concise marksheet.
# Marking scheme (marks per correct answer) changes from 3 to 5, which
# ripples into the computed "Total" row values stored on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Marking row (row 11): marks awarded per correct answer, column B ("Right")
$ws.Range("B11").Value = 5

# Total row (row 12): total marks obtained, column B ("Right")
$ws.Range("B12").Value = 115

# Total row (row 12): obtained/maximum marks text, column E ("Max")
$ws.Range("E12").Value = "115/140"
